$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated purchase/sale figures for the item in row 2
$ws.Range("C2").Value = 123.99
$ws.Range("D2").Value = 139.99

# New "Profit" column (G) with the profit expressed as a percentage
$ws.Range("G1").Value = "Profit"

$bought = $ws.Range("C2").Value()
$sold = $ws.Range("D2").Value()
$profit = [Math]::Round(($sold - $bought) / $bought * 100, 1)
$ws.Range("G2").Value = $profit

# Match column G's width to the other best-fit columns
$ws.Columns.Item(7).ColumnWidth = 5.2
